$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value  = 6.278799999999997
$ws.Range("B18").Value = 6.359199999999999
$ws.Range("B20").Value = 9.628499999999988
$ws.Range("B27").Value = 6.389400000000003
$ws.Range("B69").Value = 5.477599999999993
$ws.Range("B76").Value = 5.074200000000002
$ws.Range("B82").Value = 5.420700000000004
